# comence - order pipeline works
#
# Hire sheet: "Radio" products are renamed to "UHF"; the Hire tab becomes
# the active sheet/selection, and a handful of rows pick up their
# (re-measured) default row heights. The Bands sheet also picks up
# re-measured default row heights for its data rows and is no longer the
# active sheet.

$wb = $excel.ActiveWorkbook

$wsHire  = $wb.Worksheets.Item("Hire")
$wsBands = $wb.Worksheets.Item("Bands")

# --- Hire sheet: rename "Radio" rows (A2:A17) to "UHF" ---------------------
for ($r = 2; $r -le 17; $r++) {
    $wsHire.Cells.Item($r, 1).Value = "UHF"
}

# --- Hire sheet: row-height touch-up on rows 30-31 --------------------------
$wsHire.Rows.Item(30).RowHeight = 13.5
$wsHire.Rows.Item(31).RowHeight = 13.5

# --- Bands sheet: row-height touch-up on rows 1-14 ---------------------------
$wsBands.Rows.Item(1).RowHeight = 12.75
$wsBands.Rows.Item(2).RowHeight = 12.75
$wsBands.Rows.Item(3).RowHeight = 13.5
$wsBands.Rows.Item(4).RowHeight = 13.5
$wsBands.Rows.Item(5).RowHeight = 12.75
$wsBands.Rows.Item(6).RowHeight = 13.5
$wsBands.Rows.Item(7).RowHeight = 13.5
$wsBands.Rows.Item(8).RowHeight = 13.5
$wsBands.Rows.Item(9).RowHeight = 12.75
$wsBands.Rows.Item(10).RowHeight = 13.5
$wsBands.Rows.Item(11).RowHeight = 12.75
$wsBands.Rows.Item(12).RowHeight = 12.75
$wsBands.Rows.Item(13).RowHeight = 13.5
$wsBands.Rows.Item(14).RowHeight = 13.5

# --- Hire sheet becomes the active sheet / selection ------------------------
$wsHire.Activate()
$wsHire.Range("A2:A17").Select()
